# Add two new webcam-location rows (52 and 53) to the "location-1" sheet,
# matching the rows already present (format copied from the last existing
# data row, 51) and reusing shared strings ("LIVE, SEA, BEACH" / "Spain")
# where the new rows' text duplicates existing entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 52: Peguera, Spain webcam -----------------------------------
# Clone formatting (styles) from row 51 (the current last row) onto row 52.
$ws.Range("A51:F51").Copy() | Out-Null
$ws.Range("A52:F52").PasteSpecial(-4122) | Out-Null

$ws.Range("F52").Value = "WG6PqR9v5Uk"
$ws.Range("C52").Value = "LIVE WEBCAM Playa Palmira"
$ws.Range("B52").Value = "39.53702073445707, 2.4489938589458964"
$ws.Range("D52").Value = "Peguera"
$ws.Range("A52").Value = "LIVE, SEA, BEACH"
$ws.Range("E52").Value = "Spain"

# ---- Row 53: Fairmont The Palm, Dubai, UAE webcam ---------------------
$ws.Range("A51:F51").Copy() | Out-Null
$ws.Range("A53:F53").PasteSpecial(-4122) | Out-Null

$ws.Range("F53").Value = "7dE4IjDQJmE"
$ws.Range("E53").Value = "UAE"
$ws.Range("B53").Value = "25.108953838345435, 55.1394663849426"
$ws.Range("C53").Value = "Fairmont The Palm - The Palm Jumeirah"
$ws.Range("D53").Value = "Dubai"
$ws.Range("A53").Value = "LIVE, SEA, BEACH"

$excel.CutCopyMode = 0

# Match the author's final cursor position recorded in the saved file.
$ws.Range("C40").Select() | Out-Null
